# Generate Report for Handback
# Marks the zh-cn / de-de localization targets as handed back, fills in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for both language sheets, and widens the columns that now hold
# longer text so the report reads cleanly.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$zhCnHandbackTime = "2016-08-26 22:46:25"
$deDeHandbackTime = "2016-08-26 22:46:32"

# ---------------------------------------------------------------------------
# Overview sheet: both language status columns move from "Ready for handoff"
# to the handed-back status, and get wide enough to show the new text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet: populate "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" now that the handback round-trip finished.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3c51da9c0f270d94a1d70690b929ccd78de254c/e2e/89df6b3e-7a68-49f8-8273-30b129d815f6.md", "", "", "89df6b3e-7a68-49f8-8273-30b129d815f6.md")
$wsZhCn.Range("J2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.69cd22beefe8b94dee6ffc7afc6df686b01d2cc0.zh-cn.xlf"
$wsZhCn.Range("K2").Value = $zhCnHandbackTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3c51da9c0f270d94a1d70690b929ccd78de254c/e2e/cabe8da7-b78e-46d2-8df4-48ff53588842.md", "", "", "cabe8da7-b78e-46d2-8df4-48ff53588842.md")
$wsZhCn.Range("J3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.2e1843878eaeb4e45cb9723c90d56ef6ac4a171a.zh-cn.xlf"
$wsZhCn.Range("K3").Value = $zhCnHandbackTime

# ---------------------------------------------------------------------------
# de-de sheet: same three columns, with its own (later) handback timestamp.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3c51da9c0f270d94a1d70690b929ccd78de254c/e2e/89df6b3e-7a68-49f8-8273-30b129d815f6.md", "", "", "89df6b3e-7a68-49f8-8273-30b129d815f6.md")
$wsDeDe.Range("J2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.69cd22beefe8b94dee6ffc7afc6df686b01d2cc0.de-de.xlf"
$wsDeDe.Range("K2").Value = $deDeHandbackTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3c51da9c0f270d94a1d70690b929ccd78de254c/e2e/cabe8da7-b78e-46d2-8df4-48ff53588842.md", "", "", "cabe8da7-b78e-46d2-8df4-48ff53588842.md")
$wsDeDe.Range("J3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.2e1843878eaeb4e45cb9723c90d56ef6ac4a171a.de-de.xlf"
$wsDeDe.Range("K3").Value = $deDeHandbackTime
